# Applies the cell-value updates described in the commit:
# "Update gh-pages to output generated at 456a3b4"
# This updates the '想去人数' (F column) figures across the
# '展览' (Exhibitions), '演出' (Performances) and
# '全部类型' (All types) worksheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 5883   # was 5881
$ws.Range("F5").Value = 5883   # was 5881
$ws.Range("F7").Value = 2928   # was 2927
$ws.Range("F9").Value = 396   # was 395
$ws.Range("F13").Value = 685   # was 682
$ws.Range("F14").Value = 178   # was 176
$ws.Range("F15").Value = 4221   # was 4216
$ws.Range("F16").Value = 4221   # was 4216
$ws.Range("F22").Value = 58   # was 57
$ws.Range("F23").Value = 6338   # was 6334
$ws.Range("F24").Value = 6338   # was 6334
$ws.Range("F26").Value = 91   # was 89
$ws.Range("F29").Value = 215   # was 214
$ws.Range("F32").Value = 1614   # was 1612
$ws.Range("F34").Value = 1844   # was 1843
$ws.Range("F35").Value = 5879   # was 5873
$ws.Range("F39").Value = 75   # was 73
$ws.Range("F40").Value = 213   # was 207
$ws.Range("F41").Value = 3970   # was 3964
$ws.Range("F45").Value = 2392   # was 2391
$ws.Range("F50").Value = 280   # was 277
$ws.Range("F51").Value = 2032   # was 2031

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 185   # was 184
$ws.Range("F5").Value = 90   # was 89
$ws.Range("F10").Value = 140   # was 139

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 5883   # was 5881
$ws.Range("F5").Value = 5883   # was 5881
$ws.Range("F7").Value = 2928   # was 2927
$ws.Range("F12").Value = 185   # was 184
$ws.Range("F13").Value = 178   # was 176
$ws.Range("F14").Value = 4221   # was 4216
$ws.Range("F15").Value = 4221   # was 4216
$ws.Range("F21").Value = 58   # was 57
$ws.Range("F22").Value = 6338   # was 6334
$ws.Range("F23").Value = 6338   # was 6334
$ws.Range("F25").Value = 91   # was 89
$ws.Range("F27").Value = 215   # was 214
$ws.Range("F28").Value = 90   # was 89
$ws.Range("F30").Value = 1614   # was 1612
$ws.Range("F33").Value = 1844   # was 1843
$ws.Range("F35").Value = 5879   # was 5873
$ws.Range("F39").Value = 75   # was 73
$ws.Range("F40").Value = 3970   # was 3964
$ws.Range("F45").Value = 2392   # was 2391
$ws.Range("F50").Value = 280   # was 277
$ws.Range("F51").Value = 140   # was 139

$wb.Save()
